$wb = $excel.ActiveWorkbook

# New Gurobi experiment results (rerun with academic license)
$data = @{}
$data["N1_D40"] = @{
    2 = @(0.422, 16.5)
    3 = @(0.017, 15.39)
    4 = @(0.017, 15.27)
    5 = @(0.017, 15.41)
    6 = @(0.017, 15.46)
    7 = @(0.018, 15.41)
    8 = @(0.018, 15.34)
    9 = @(0.018, 15.46)
    10 = @(0.019, 15.46)
    11 = @(0.017, 15.48)
    12 = @(0.05800000000000001, 15.518)
}
$data["N1_D60"] = @{
    2 = @(0.026, 14.91)
    3 = @(0.026, 14.92)
    4 = @(0.027, 14.92)
    5 = @(0.026, 14.87)
    6 = @(0.026, 14.94)
    7 = @(0.026, 14.73)
    8 = @(0.027, 14.86)
    9 = @(0.026, 14.8)
    10 = @(0.026, 14.79)
    11 = @(0.026, 14.84)
    12 = @(0.0262, 14.858)
}
$data["N1_D80"] = @{
    2 = @(0.04, 18.41)
    3 = @(0.04, 18.47)
    4 = @(0.039, 18.34)
    5 = @(0.039, 18.46)
    6 = @(0.04, 18.34)
    7 = @(0.039, 18.42)
    8 = @(0.039, 18.39)
    9 = @(0.039, 18.48)
    10 = @(0.039, 18.41)
    11 = @(0.04, 18.42)
    12 = @(0.0394, 18.414)
}
$data["N1_D100"] = @{
    2 = @(0.061, 16.48)
    3 = @(0.058, 16.44)
    4 = @(0.058, 16.33)
    5 = @(0.058, 16.33)
    6 = @(0.058, 16.47)
    7 = @(0.058, 16.53)
    8 = @(0.058, 16.47)
    9 = @(0.057, 16.44)
    10 = @(0.058, 16.33)
    11 = @(0.058, 16.46)
    12 = @(0.05820000000000001, 16.428)
}

foreach ($sheetName in $data.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $data[$sheetName]
    foreach ($rowNum in $rows.Keys) {
        $vals = $rows[$rowNum]
        $ws.Cells.Item($rowNum, 5).Value = $vals[0]
        $ws.Cells.Item($rowNum, 6).Value = $vals[1]
    }
}
